$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value that must be stored as TEXT even when it looks like
# a number/date (this sheet keeps every column except I/M/N as text).
# Briefly switching the cell to the "Text" number format forces Excel to
# keep the literal string instead of auto-converting it to a number/date,
# then restoring the "Normal" cell style returns the cell to its original
# (unstyled) appearance, matching the rest of the sheet.
function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# --- Update existing row 50 ---
Set-TextValue $ws.Range("A50") "7023"
Set-TextValue $ws.Range("C50") "ZABALA 3573"
$ws.Range("M50").Value = -58.457522
$ws.Range("N50").Value = -34.579414

# --- Append new row 75 ---
Set-TextValue $ws.Range("A75") "7021"
Set-TextValue $ws.Range("B75") "8/25/2025"
Set-TextValue $ws.Range("C75") "VERA 445"
Set-TextValue $ws.Range("D75") "15"
Set-TextValue $ws.Range("E75") "809155622"
Set-TextValue $ws.Range("F75") "NEW"
Set-TextValue $ws.Range("G75") "Pendiente"
Set-TextValue $ws.Range("H75") "Picada"
$ws.Range("I75").Value = 1
Set-TextValue $ws.Range("J75") "Cambio"
Set-TextValue $ws.Range("K75") "Sin equipos"
Set-TextValue $ws.Range("L75") "Pasante"
$ws.Range("M75").Value = -58.437181
$ws.Range("N75").Value = -34.5995
Set-TextValue $ws.Range("O75") "Palermo"
Set-TextValue $ws.Range("P75") "Capital Sur"
